$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3932.5938
$ws.Range("I74").Value = 3861.7693
$ws.Range("J74").Value = 3981.0527
$ws.Range("K74").Value = 3861.7693
$ws.Range("L74").Value = 3981.0527
$ws.Range("M74").Value = -2925.7693
$ws.Range("N74").Value = -5853.0527
$ws.Range("H77").Value = 3932.5938
$ws.Range("I77").Value = 3861.7693
$ws.Range("J77").Value = 3981.0527
$ws.Range("K77").Value = 19308.8465
$ws.Range("L77").Value = 19905.2635
$ws.Range("M77").Value = -14628.8465
$ws.Range("N77").Value = -29265.2635
$ws.Range("H100").Value = 1333.3334
$ws.Range("I100").Value = 1380
$ws.Range("J100").Value = 1100
$ws.Range("K100").Value = 1380
$ws.Range("L100").Value = 1100
$ws.Range("M100").Value = -839
$ws.Range("N100").Value = -2182
$ws.Range("H121").Value = 1529.9
$ws.Range("I121").Value = 695
$ws.Range("J121").Value = 1622.6666
$ws.Range("K121").Value = 2085
$ws.Range("L121").Value = 4867.9998
$ws.Range("M121").Value = -338
$ws.Range("N121").Value = -8361.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1397.9565
$ws.Range("I2").Value = 1354.9048
$ws.Range("J2").Value = 1850
$ws.Range("K2").Value = 1354.9048
$ws.Range("L2").Value = 1850
$ws.Range("M2").Value = -1241.9048
$ws.Range("N2").Value = -2076
$ws.Range("H4").Value = 340.2
$ws.Range("I4").Value = 300.25
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 300.25
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = -184.25
$ws.Range("N4").Value = -732
$ws.Range("H5").Value = 26315864
$ws.Range("I5").Value = 26315864
$ws.Range("K5").Value = 26315864
$ws.Range("M5").Value = -26315752
$ws.Range("H32").Value = 14817.238
$ws.Range("I32").Value = 15289.134
$ws.Range("J32").Value = 12094.77
$ws.Range("K32").Value = 15289.134
$ws.Range("L32").Value = 12094.77
$ws.Range("M32").Value = -15002.134
$ws.Range("N32").Value = -12668.77
$ws.Range("H45").Value = 1763.0312
$ws.Range("I45").Value = 1842.6072
$ws.Range("J45").Value = 1206
$ws.Range("K45").Value = 1842.6072
$ws.Range("L45").Value = 1206
$ws.Range("M45").Value = -1465.6072
$ws.Range("N45").Value = -1960
$ws.Range("H74").Value = 10501262
$ws.Range("I74").Value = 13211879
$ws.Range("J74").Value = 200916
$ws.Range("K74").Value = 13211879
$ws.Range("L74").Value = 200916
$ws.Range("M74").Value = -13211005
$ws.Range("N74").Value = -202664
$ws.Range("H77").Value = 10501262
$ws.Range("I77").Value = 13211879
$ws.Range("J77").Value = 200916
$ws.Range("K77").Value = 66059395
$ws.Range("L77").Value = 1004580
$ws.Range("M77").Value = -66055027
$ws.Range("N77").Value = -1013316
$ws.Range("H102").Value = 1579
$ws.Range("I102").Value = 1390.5454
$ws.Range("J102").Value = 2097.25
$ws.Range("K102").Value = 1390.5454
$ws.Range("L102").Value = 2097.25
$ws.Range("M102").Value = 231.4546
$ws.Range("N102").Value = -5341.25
$ws.Range("H116").Value = 1397.9565
$ws.Range("I116").Value = 1354.9048
$ws.Range("J116").Value = 1850
$ws.Range("K116").Value = 1354.9048
$ws.Range("L116").Value = 1850
$ws.Range("M116").Value = 939.0952
$ws.Range("N116").Value = -6438

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1397.9565
$ws.Range("I3").Value = 1354.9048
$ws.Range("J3").Value = 1850
$ws.Range("K3").Value = 1354.9048
$ws.Range("L3").Value = 1850
$ws.Range("M3").Value = -1240.9048
$ws.Range("N3").Value = -2078
$ws.Range("H4").Value = 26315864
$ws.Range("I4").Value = 26315864
$ws.Range("K4").Value = 26315864
$ws.Range("M4").Value = -26315749
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H86").Value = 19586.154
$ws.Range("I86").Value = 30025.75
$ws.Range("K86").Value = 30025.75
$ws.Range("M86").Value = -28902.75
$ws.Range("H89").Value = 19586.154
$ws.Range("I89").Value = 30025.75
$ws.Range("K89").Value = 150128.75
$ws.Range("M89").Value = -144512.75
$ws.Range("H103").Value = 3742.5
$ws.Range("J103").Value = 3742.5
$ws.Range("L103").Value = 3742.5
$ws.Range("N103").Value = -6086.5
$ws.Range("H105").Value = 189773.75
$ws.Range("I105").Value = 216462.86
$ws.Range("K105").Value = 216462.86
$ws.Range("M105").Value = -214715.86
$ws.Range("H107").Value = 2719.6
$ws.Range("I107").Value = 2650.1667
$ws.Range("K107").Value = 2650.1667
$ws.Range("M107").Value = -730.1667000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 3353333.2
$ws.Range("J4").Value = 3353333.2
$ws.Range("L4").Value = 3353333.2
$ws.Range("N4").Value = -3353557.2
$ws.Range("H7").Value = 125.7619
$ws.Range("I7").Value = 36.666668
$ws.Range("J7").Value = 348.5
$ws.Range("K7").Value = 36.666668
$ws.Range("L7").Value = 348.5
$ws.Range("M7").Value = 76.333332
$ws.Range("N7").Value = -574.5
$ws.Range("H16").Value = 1169.1765
$ws.Range("I16").Value = 711.8182
$ws.Range("K16").Value = 711.8182
$ws.Range("M16").Value = -424.8182
$ws.Range("H20").Value = 51000
$ws.Range("J20").Value = 51000
$ws.Range("L20").Value = 51000
$ws.Range("N20").Value = -51472
$ws.Range("H30").Value = 51000
$ws.Range("J30").Value = 51000
$ws.Range("L30").Value = 51000
$ws.Range("N30").Value = -51182
$ws.Range("H31").Value = 3456.2334
$ws.Range("I31").Value = 1819.381
$ws.Range("J31").Value = 7275.5557
$ws.Range("K31").Value = 1819.381
$ws.Range("L31").Value = 7275.5557
$ws.Range("M31").Value = -1524.381
$ws.Range("N31").Value = -7865.5557
$ws.Range("H34").Value = 3456.2334
$ws.Range("I34").Value = 1819.381
$ws.Range("J34").Value = 7275.5557
$ws.Range("K34").Value = 1819.381
$ws.Range("L34").Value = 7275.5557
$ws.Range("M34").Value = -1617.381
$ws.Range("N34").Value = -7679.5557
$ws.Range("H113").Value = 1169.1765
$ws.Range("I113").Value = 711.8182
$ws.Range("K113").Value = 711.8182
$ws.Range("M113").Value = 1458.1818
$ws.Range("H127").Value = 32857.145
$ws.Range("J127").Value = 32857.145
$ws.Range("L127").Value = 32857.145
$ws.Range("N127").Value = -42777.145
$ws.Range("H128").Value = 51000
$ws.Range("J128").Value = 51000
$ws.Range("L128").Value = 51000
$ws.Range("N128").Value = -60960

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 41083.46
$ws.Range("I132").Value = 2725.125
$ws.Range("J132").Value = 102456.8
$ws.Range("K132").Value = 8175.375
$ws.Range("L132").Value = 307370.4
$ws.Range("M132").Value = -5645.375
$ws.Range("N132").Value = -312430.4
$ws.Range("H4").Value = 7648818
$ws.Range("J4").Value = 6668593.5
$ws.Range("L4").Value = 20005780.5
$ws.Range("N4").Value = -20006004.5
$ws.Range("H5").Value = 710.8
$ws.Range("I5").Value = 605.63635
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 1816.90905
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = -1704.90905
$ws.Range("N5").Value = -3224

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 710.8
$ws.Range("I135").Value = 605.63635
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 5450.72715
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -2915.72715
$ws.Range("N135").Value = -14070
$ws.Range("H113").Value = 1237.05
$ws.Range("I113").Value = 1046.3125
$ws.Range("K113").Value = 1046.3125
$ws.Range("M113").Value = 1123.6875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1400
$ws.Range("I93").Value = 1400
$ws.Range("K93").Value = 1400
$ws.Range("M93").Value = -152
